$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- "external" command list (#system column I): add tail(id,file) ---
$ws.Range("I5").Value = "tail(id,file)"

# --- "web" command list (#system column Y): insert
#     assertTextNotContains(locator,text) alphabetically before
#     assertTextNotPresent(text), shifting Y39:Y127 down to Y40:Y128 ---
for ($i = 127; $i -ge 39; $i--) {
    $src = $ws.Range("Y" + $i).Value2
    $ws.Range("Y" + ($i + 1)).Value = $src
}
$ws.Range("Y39").Value = "assertTextNotContains(locator,text)"

# --- update the defined names so their ranges include the new rows ---
$wb.Names.Item("external").RefersTo = "='#system'!`$I`$2:`$I`$5"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$128"
